$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B; this shifts old B,C,D -> C,D,E
$ws.Columns.Item(2).Insert()

# The Insert() above copies column A's formatting into new column B;
# the segment-name labels that will live there are unstyled, so clear it.
$ws.Range("B2:B20").ClearFormats()

# Move the segment-name labels currently in column A (rows 2-20) to column B
for ($r = 2; $r -le 20; $r++) {
    $name = $ws.Cells.Item($r, 1).Value()
    $ws.Cells.Item($r, 2).Value = $name
    $ws.Cells.Item($r, 1).Value = $r - 2
}

# Header for the new column - copy formatting from neighboring header cell, then set text
$ws.Cells.Item(1, 3).Copy()
$ws.Cells.Item(1, 2).PasteSpecial(-4122)
$ws.Cells.Item(1, 2).Value = "segments"

Write-Output "done"
